$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new "Dee" dialogue line (pool of water) ---
$ws.Range("A2").Value = "Dee"
$ws.Range("B2").Value = "There’s a pool of water on the ground outside the door, now still not dried."
$ws.Range("C2").Value = "Dee-Thinking2"
$ws.Range("E2").Value = "Water"
$ws.Range("F2").Value = "Suspicious"
$ws.Range("L2").Value = "Dee-Regular"

# --- Row 3: new "Dee" dialogue line (scene remains) ---
$ws.Range("A3").Value = "Dee"
$ws.Range("B3").Value = "Other than that, the scene remains just as it was last night."
$ws.Range("C3").Value = "Dee-Thinking2"
$ws.Range("E3").Value = "Water"

# Rows 2 & 3 grow taller to fit the wrapped text
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34

# --- Row 4 ---
$ws.Range("A4").Value = "Investigate2"
$ws.Range("C4").Value = "Water"

# --- Row 5 ---
$ws.Range("A5").ClearContents()
$ws.Range("B5").Value = "Hand"
$ws.Range("C5").Value = "Hand"

# --- Row 6 ---
$ws.Range("B6").Value = "Blood"
$ws.Range("C6").Value = "Blood"
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

# --- Row 7 ---
$ws.Range("B7").Value = "End Investigation"
$ws.Range("C7").Value = "StoryScript14"

# --- Row 8 no longer exists ---
$ws.Rows.Item(8).Delete()

$ws.Range("B10").Select() | Out-Null
